$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 0.4818104178965973
$ws.Range("C9").Value = 0.4818104178965973
$ws.Range("D9").Value = 0.4086707003526728
$ws.Range("E9").Value = 0.6392735723871845
$ws.Range("F9").Value = 0.4602556961850388
$ws.Range("G9").Value = 6

$ws.Range("B10").Value = 0.2297780799482224
$ws.Range("C10").Value = 0.2410916961729521
$ws.Range("D10").Value = 0.09195223546803079
$ws.Range("E10").Value = 0.303236270040427
$ws.Range("F10").Value = 0.2423456295562365
$ws.Range("G10").Value = 3

$ws.Range("B11").Value = 0.2888921154092369
$ws.Range("C11").Value = 0.2888921154092369
$ws.Range("D11").Value = 0.08345865434562384
$ws.Range("E11").Value = 0.2888921154092369
$ws.Range("F11").ClearContents()
$ws.Range("G11").Value = 1
